# Remove unused interventions from Philippines sheet
#
# The "program_perc_xpertacf_indigenous" time-variant parameter (row 9 of
# the "time_variants" sheet) is unused, so the whole row is deleted. This
# shifts every row below it up by one and - because that shared string then
# has no remaining references anywhere in the workbook - it also drops out
# of the shared-strings table on save, which is exactly what the target
# diff shows (uniqueCount 69 -> 68, and every string index above it
# decremented by one across all sheets).
#
# The edit also leaves the "time_variants" sheet as the active tab/sheet,
# with a fresh selection, and tidies up the "constants" sheet's selection.

$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Sheets.Item("constants")
$wsTimeVariants = $wb.Sheets.Item("time_variants")

# Delete the entire "program_perc_xpertacf_indigenous" row.
$wsTimeVariants.Rows(9).Delete()

# "constants" is no longer the selected tab; park its selection at B2.
$wsConstants.Range("B2").Select()

# "time_variants" becomes the active sheet/tab, selected at A8.
$wsTimeVariants.Activate()
$wsTimeVariants.Range("A8").Select()
